$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.641.48"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.890.81"
$ws.Range("E3").Value = "  +0.06%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.003"
$ws.Range("E4").Value = "  +0.07%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "239.13"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  +0.08%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4907"
$ws.Range("E7").Value = "  +0.46%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2938"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06702"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "1.900.11"
$ws.Range("E10").Value = "  +0.59%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "16.99"
$ws.Range("E11").Value = "  +0.09%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07338"
$ws.Range("E12").Value = "  +1.56%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "5.173"
$ws.Range("E13").Value = "  +3.31%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "88.05"
$ws.Range("E14").Value = "  -1.82%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6666"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "30.608.86"
$ws.Range("E16").Value = "  +0.07%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000007875"
$ws.Range("E17").Value = "  -0.52%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "13.43"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").Value = "2.156.90"
$ws.Range("E20").Value = "  +0.81%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.356"
$ws.Range("E21").Value = "  +12.13%  "
$ws.Range("E22").Value = "  -0.07%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "189.88"
$ws.Range("E23").Value = "  +0.04%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.203"
$ws.Range("E24").Value = "  +2.75%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.529"
$ws.Range("E25").Value = "  +2.29%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "161.66"
$ws.Range("E26").Value = "  +3.43%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.47"
$ws.Range("E27").Value = "  -0.87%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.930"
$ws.Range("E28").Value = "  +2.95%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.466"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.404"
$ws.Range("E30").Value = "  +3.43%  "
$ws.Range("E31").Value = "  +1.46%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.040"
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("E33").Value = "  -0.72%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7426"
$ws.Range("E34").Value = "  +0.81%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.100"
$ws.Range("E35").Value = "  +0.99%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.731"
$ws.Range("E36").Value = "  -1.04%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.01820"
$ws.Range("E37").Value = "  -0.08%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.700"
$ws.Range("E38").Value = "  +0.81%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.9165"
$ws.Range("E39").Value = "  -0.57%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.066"
$ws.Range("E40").Value = "  -0.41%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "75.32"
$ws.Range("E41").Value = "  +31.12%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.4416"
$ws.Range("E42").Value = "  +0.90%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.916"
$ws.Range("E43").Value = "  +4.31%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "106.32"
$ws.Range("E44").Value = "  +1.49%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.9937"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("E47").Value = "  +1.52%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "35.41"
$ws.Range("E48").Value = "  +5.08%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.997"
$ws.Range("E49").Value = "  +2.72%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.05843"
$ws.Range("E50").Value = "  -0.14%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.3961"
$ws.Range("E51").Value = "  +0.73%  "